# Update gh-pages generated output (南宁-漫展信息.xlsx)
# Sheet "展览" (exhibitions) and "全部类型" (all types) both list the same
# events; the scraped stats (interested-count / min price) changed on refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 305
$ws1.Range("F4").Value = 1227
$ws1.Range("G4").Value = 57
$ws1.Range("F5").Value = 613

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 305
$ws4.Range("F4").Value = 1227
$ws4.Range("G4").Value = 57
$ws4.Range("F6").Value = 613
